# Update cryptocurrency listing with the latest prices/volumes/ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "29.332.26"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.844.19"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Formula = '=TEXT(0.9968,"0.0000")'
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Formula = '=TEXT(239.83,"000.00")'
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Formula = '=TEXT(0.6279,"0.0000")'
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Formula = '=TEXT(0.9987,"0.0000")'
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Formula = '=TEXT(0.07496,"0.00000")'
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Formula = '=TEXT(0.2899,"0.0000")'
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Formula = '=TEXT(24.42,"00.00")'
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Formula = '=TEXT(0.07730,"0.00000")'
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.843.62"
$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Formula = '=TEXT(4.987,"0.000")'
$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Formula = '=TEXT(0.6800,"0.0000")'
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Formula = '=TEXT(0.00001049,"0.00000000")'
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Formula = '=TEXT(82.02,"00.00")'
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Formula = '=TEXT(6.180,"0.000")'
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.394.72"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Formula = '=TEXT(229.04,"000.00")'
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Formula = '=TEXT(12.32,"00.00")'
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Formula = '=TEXT(0.9984,"0.0000")'
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Formula = '=TEXT(7.495,"0.000")'
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Formula = '=TEXT(0.9988,"0.0000")'
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Formula = '=TEXT(158.50,"000.00")'
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Formula = '=TEXT(8.422,"0.000")'
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Formula = '=TEXT(0.1371,"0.0000")'
$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Formula = '=TEXT(17.52,"00.00")'
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Formula = '=TEXT(0.06395,"0.00000")'
$ws.Range("E28").Value = "  +14.16%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Formula = '=TEXT(1.407,"0.000")'
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Formula = '=TEXT(1.479,"0.000")'
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Formula = '=TEXT(4.090,"0.000")'
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Formula = '=TEXT(4.090,"0.000")'
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Formula = '=TEXT(1.832,"0.000")'
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Formula = '=TEXT(1.140,"0.000")'
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Formula = '=TEXT(0.6970,"0.0000")'
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Formula = '=TEXT(2.578,"0.000")'
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.268.05"
$ws.Range("E37").Value = "  +3.38%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Formula = '=TEXT(2.843,"0.000")'
$ws.Range("E38").Value = "  +4.60%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Formula = '=TEXT(0.01833,"0.00000")'
$ws.Range("E39").Value = "  +1.67%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Formula = '=TEXT(6.722,"0.000")'
$ws.Range("E40").Value = "  +5.73%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Formula = '=TEXT(0.9148,"0.0000")'
$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Formula = '=TEXT(0.9982,"0.0000")'
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.005.81"
$ws.Range("E43").Value = "  -18.49%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Formula = '=TEXT(101.16,"000.00")'
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Formula = '=TEXT(66.14,"00.00")'
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Formula = '=TEXT(7.080,"0.000")'
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Formula = '=TEXT(1.722,"0.000")'
$ws.Range("E47").Value = "  +2.40%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Formula = '=TEXT(0.1165,"0.0000")'
$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Formula = '=TEXT(0.00000000116,"0.00000000000")'
$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Formula = '=TEXT(0.3960,"0.0000")'
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Formula = '=TEXT(8.977,"0.000")'
$ws.Range("E51").Value = "  -0.14%  "

# Flatten the TEXT() helper formulas above into plain text values
# (xlPasteValues = -4163) so the cells keep their original text type
# instead of staying as live formulas.
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$excel.CutCopyMode = 0

